$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Utilizacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de:</w:t></w:r></w:p>')
$d.Paragraphs.Item(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Creator</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$d.Paragraphs.Item(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 5.7.1</w:t></w:r></w:p>')
$d.Paragraphs.Item(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Mingw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 32bits</w:t></w:r></w:p>')
$d.Paragraphs.Item(5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Glm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 0.9.8</w:t></w:r></w:p>')
$d.Paragraphs.Item(8).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Creación de repositorio para control de versiones, utilizando en este caso </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$d.Paragraphs.Item(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Creacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> del proyecto en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>creator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$d.Paragraphs.Item(10).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Creacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de la clase </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GLWidget</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$d.Paragraphs.Item(11).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Configuración básica para la utilización de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opengl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$d.Paragraphs.Item(12).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Creacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de ventana para visualización de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opengl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y zona para parámetros.</w:t></w:r></w:p>')
$d.Paragraphs.Item(13).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Creacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perspective</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” a partir de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>glFrustrum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> puesto que no se puede utilizar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>glPerspective</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> utiliza </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opengl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3.X y a partir de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>opengl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3 no se puede).</w:t></w:r></w:p>')
$d.Paragraphs.Item(14).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Utilizacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de la librería GLM para utilizar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LookAt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que creara la matriz View y añadiremos a la pila con </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>glLoadMatrix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>')
$d.Paragraphs.Item(15).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">// se podría utilizar la librería </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>glm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para usar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perspective</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> puesto que esta la tiene implementada.</w:t></w:r></w:p>')
$d.Paragraphs.Item(17).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Activar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>strong</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>focus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> mediante “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setFocusPolicy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="800080"/></w:rPr><w:t>Qt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>::</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="800080"/></w:rPr><w:t>StrongFocus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">);” para que el widget </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>glwidget</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pueda recibir eventos de teclado y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>raton</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$d.Paragraphs.Item(18).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Introduccion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>formulas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de movimiento sobre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>robotica</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

# Restructure the final block: paragraphs 20-22 (1-indexed) -> 4 paragraphs
$pStart = $d.Paragraphs.Item(20)
$pEnd = $d.Paragraphs.Item(22)
$rBlock = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rBlock.InsertXML('<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>A</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Widgets para la introducción de parámetros (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wheelSpeed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wheelRadius</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wheelSeparation,sensorSeparation,robot</w:t></w:r><w:r><w:t>Diameter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>iu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Paso de parámetros entre widgets.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p>')